# Updated to Carl's newest scraper: the source sheet now only ships the
# 9 columns that matter (Varumärke .. RSK-nummer); drop every other
# legacy/import column and the old autofilter/sort state that referenced
# the wide S-column range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn off the worksheet AutoFilter (also drops the <sortState> it carried).
$ws.AutoFilterMode = $false

# Remove the now-unused columns, right-to-left so earlier deletions don't
# shift the column letters we still need to remove.
# Columns to drop (by their original letter): Q, O, M, L, K, J, I, H, G, A
# Columns kept, in order: B, C, D, E, F, N, P, R, S -> become A..I
$colsToDelete = @("Q", "O", "M", "L", "K", "J", "I", "H", "G", "A")
foreach ($col in $colsToDelete) {
    $ws.Columns($col + ":" + $col).Delete()
}

# Widen the new "SSG-notering" column (old P, now G) to fit its content.
$ws.Columns("G:G").ColumnWidth = 28.666666666666668

# Reset the view: scroll back to column A and select the new last column.
[void]$ws.Range("H1:H1048576").Select()

# Keep the _FilterDatabase defined name in sync with the new, narrower range.
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$I`$21217"
